$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.665.30'
$ws.Range('E2').Value = '  +0.39%  '

$ws.Range('D3').Value = '1.844.07'

$ws.Range('E4').Value = '  -0.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.19'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.27%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.18%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4278'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.96%  '

$ws.Range('E8').Value = '  +0.16%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07323'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.96%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8777'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.27%  '

$ws.Range('E11').Value = '  +0.79%  '

$ws.Range('D12').Value = '1.861.56'
$ws.Range('E12').Value = '  -2.23%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.349'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.30%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.518'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.60%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06954'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.99%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.004'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.11%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '79.52'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.77%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008995'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.23%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.002'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.10%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.38'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.76%  '

$ws.Range('D21').Value = '27.740.65'
$ws.Range('E21').Value = '  +0.70%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.984'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.19%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.31'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.92%  '

$ws.Range('D24').Value = '2.081.05'
$ws.Range('E24').Value = '  +1.46%  '

$ws.Range('E25').Value = '  -2.47%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '155.75'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.32%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.55'
$ws.Range('D27').Style = 'Normal'

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '119.65'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.39%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.227'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.55%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.874'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.37%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08886'
$ws.Range('D31').Style = 'Normal'

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7534'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.59%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.528'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.42%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.952'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.23%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.119'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.15%  '

$ws.Range('E36').Value = '  -0.07%  '

$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.108'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.09%  '

$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05432'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.84%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01935'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.79%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.832'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.32%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1667'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.30%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5070'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.50%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.600'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.18%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.374'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.50%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.06536'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.27%  '

$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '106.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.01%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.34'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.42%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4649'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.73%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.001'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.22%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.636'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.61%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '64.59'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.21%  '
